$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: update the values in the first three rows ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Step 2: insert 10 new rows right after row 3 (i.e. before the
#     row that currently follows it), each holding a single value ---
$newValues = @("106", "0.00003", "0.00018", "0.00006", "0.00004", `
                "0.00004", "0.00008", "0.00018", "0.00496", "100.0")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $refRow = $t.Rows.Item(4)
    $newRow = $t.Rows.Add($refRow)
}
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowIndex = 4 + $i
    $t.Rows.Item($rowIndex).Cells.Item(1).Range.Text = $newValues[$i]
}

# --- Step 3: collapse the last three multi-value (tab-separated) rows
#     down to a single value each. These were rows 34/35/36 before the
#     10-row insertion, so they are now shifted by +10. ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.99"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "90"
